$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add new columns C (product_pic) and D (productPrice)
$ws.Range("C1").Value = "product_pic"
$ws.Range("D1").Value = "productPrice"

# Row 2: fill in C and D for the existing Sanyo S23 product row
$ws.Range("C2").Value = "/static/img/product-01.jpg"
$ws.Range("D2").Value = 12

# Rows 3-13: new product rows Sanyo S24 .. Sanyo S34
$names = @("Sanyo S24", "Sanyo S25", "Sanyo S26", "Sanyo S27", "Sanyo S28", "Sanyo S29", "Sanyo S30", "Sanyo S31", "Sanyo S32", "Sanyo S33", "Sanyo S34")

$row = 3
$price = 13
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = "Phone"
    $ws.Cells.Item($row, 3).Value = "/static/img/product-01.jpg"
    $ws.Cells.Item($row, 4).Value = $price
    $row = $row + 1
    $price = $price + 1
}

# Match the selection recorded in the saved file
$ws.Range("I9").Select()
